# software pricing auto fill
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update part number string in A2
$ws.Range("A2").Value = "2.3.01.01.10597"

# Update pricing values across B2:H2
$ws.Range("B2").Value = 9000
$ws.Range("C2").Value = 10211
$ws.Range("D2").Value = 11604
$ws.Range("E2").Value = 13091
$ws.Range("F2").Value = 14587
$ws.Range("G2").Value = 15709
$ws.Range("H2").Value = 39273
